$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New attendance column for 2025-04-07.
# Force text so the ISO-looking date string isn't auto-converted to a date serial,
# then drop back to the default "Normal" style (matching the plain, unstyled
# header/value cells used for this new column).
$ws.Cells.Item(1, 4).NumberFormat = "@"
$ws.Cells.Item(1, 4).Value = "2025-04-07"
$ws.Cells.Item(1, 4).Style = "Normal"

$ws.Cells.Item(2, 4).Value = "P"
$ws.Cells.Item(2, 4).Style = "Normal"
